# Add data for 2021-11-25: updates the "through" date in the sheet name and
# header label, refreshes the current-November running total, and bumps
# several historical "November <year>" column counts that moved for
# neighborhoods whose data was backfilled for 2021-11-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab and update the column-B header label to reflect
# the new "through" date.
$ws.Name = "Through 2021-11-17"
$ws.Range("B1").Value = "November 2021 (through November 17)"

# Current "November 2021 (through November 17)" column (B) updates.
$ws.Range("B9").Value = 1
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 3
$ws.Range("B20").Value = 3
$ws.Range("B31").Value = 3

# Historical "November <year>" column updates.
$ws.Range("M3").Value = 10
$ws.Range("AT3").Value = 5
$ws.Range("BE3").Value = 4
$ws.Range("AI4").Value = 5
$ws.Range("AT4").Value = 10
$ws.Range("AI5").Value = 2
$ws.Range("M7").Value = 7
$ws.Range("BE12").Value = 3
$ws.Range("M30").Value = 3
$ws.Range("AT45").Value = 5
$ws.Range("X48").Value = 3
$ws.Range("M62").Value = 1
$ws.Range("M68").Value = 3
$ws.Range("M84").Value = 1
$ws.Range("X94").Value = 1
